# replace space in group column names
# The "group_cooking" repeat-group sheet had header names with spaces
# ("Cooking Equipment", "Years Owned"); rename them to use underscores
# instead, matching the XLSForm/ODK group-name convention.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("group_cooking")

$ws2.Range("A1").Value = "Cooking_Equipment"
$ws2.Range("B1").Value = "Years_Owned"

# The sheet was also brought to the front / made the active tab, with the
# selection parked on B2.
$ws2.Activate()
$ws2.Range("B2").Select()
